$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 29907
$ws.Range("J75").Value = 29907
$ws.Range("L75").Value = 29907
$ws.Range("N75").Value = -31779

$ws.Range("H78").Value = 29907
$ws.Range("J78").Value = 29907
$ws.Range("L78").Value = 89721
$ws.Range("N78").Value = -99081

$ws.Range("H88").Value = 3696.6428
$ws.Range("I88").Value = 989.8
$ws.Range("J88").Value = 5200.4443
$ws.Range("K88").Value = 989.8
$ws.Range("L88").Value = 5200.4443
$ws.Range("M88").Value = -583.8
$ws.Range("N88").Value = -6012.4443

$ws.Range("H91").Value = 3696.6428
$ws.Range("I91").Value = 989.8
$ws.Range("J91").Value = 5200.4443
$ws.Range("K91").Value = 989.8
$ws.Range("L91").Value = 5200.4443
$ws.Range("M91").Value = 414.2
$ws.Range("N91").Value = -8008.4443

$ws.Range("H129").Value = 972.6667
$ws.Range("J129").Value = 1079.9736
$ws.Range("L129").Value = 3239.9208
$ws.Range("N129").Value = -13239.9208

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 92416.27
$ws.Range("I45").Value = 125999.875
$ws.Range("J45").Value = 2860
$ws.Range("K45").Value = 125999.875
$ws.Range("L45").Value = 2860
$ws.Range("M45").Value = -125622.875
$ws.Range("N45").Value = -3614

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 29404.75
$ws.Range("I86").Value = 42319.742
$ws.Range("K86").Value = 42319.742
$ws.Range("M86").Value = -41196.742

$ws.Range("H89").Value = 29404.75
$ws.Range("I89").Value = 42319.742
$ws.Range("K89").Value = 211598.71
$ws.Range("M89").Value = -205982.71

$ws.Range("H134").Value = 2973.9736
$ws.Range("J134").Value = 3049
$ws.Range("L134").Value = 9147
$ws.Range("N134").Value = -14217

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29517.807
$ws.Range("I31").Value = 1121.1765
$ws.Range("J31").Value = 41586.375
$ws.Range("K31").Value = 1121.1765
$ws.Range("L31").Value = 41586.375
$ws.Range("M31").Value = -826.1765
$ws.Range("N31").Value = -42176.375

$ws.Range("H34").Value = 29517.807
$ws.Range("I34").Value = 1121.1765
$ws.Range("J34").Value = 41586.375
$ws.Range("K34").Value = 1121.1765
$ws.Range("L34").Value = 41586.375
$ws.Range("M34").Value = -919.1765
$ws.Range("N34").Value = -41990.375

$ws.Range("H132").Value = 38464828
$ws.Range("I132").Value = 32261268
$ws.Range("J132").Value = 62503624
$ws.Range("K132").Value = 96783804
$ws.Range("L132").Value = 187510872
$ws.Range("M132").Value = -96781274
$ws.Range("N132").Value = -187515932

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1675.6875
$ws.Range("I5").Value = 836.5714
$ws.Range("J5").Value = 2328.3333
$ws.Range("K5").Value = 2509.7142
$ws.Range("L5").Value = 6984.999899999999
$ws.Range("M5").Value = -2397.7142
$ws.Range("N5").Value = -7208.999899999999

$ws.Range("H68").Value = 500
$ws.Range("I68").Value = 500
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 1500
$ws.Range("L68").Value = 0
$ws.Range("N68").Value = $null
$ws.Range("M68").Value = -689

$ws.Range("H71").Value = 500
$ws.Range("I71").Value = 500
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 4500
$ws.Range("L71").Value = 0
$ws.Range("N71").Value = $null
$ws.Range("M71").Value = -444

$ws.Range("H131").Value = 852.51
$ws.Range("I131").Value = 601.6667
$ws.Range("J131").Value = 868.5213
$ws.Range("K131").Value = 1805.0001
$ws.Range("L131").Value = 2605.5639
$ws.Range("M131").Value = 3234.9999
$ws.Range("N131").Value = -12685.5639

$ws.Range("H132").Value = 2175
$ws.Range("I132").Value = 1000
$ws.Range("J132").Value = 2248.4375
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 20235.9375
$ws.Range("M132").Value = -6470
$ws.Range("N132").Value = -25295.9375

$ws.Range("H135").Value = 1675.6875
$ws.Range("I135").Value = 836.5714
$ws.Range("J135").Value = 2328.3333
$ws.Range("K135").Value = 7529.1426
$ws.Range("L135").Value = 20954.9997
$ws.Range("M135").Value = -4994.1426
$ws.Range("N135").Value = -26024.9997

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 166840030
$ws.Range("I80").Value = 200207600
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 200207600
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -200206602
$ws.Range("N80").Value = -4196

$ws.Range("H83").Value = 166840030
$ws.Range("I83").Value = 200207600
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 1001038000
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -1001033008
$ws.Range("N83").Value = -20984

$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = $null

$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null

$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = $null

$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null

$ws.Range("H132").Value = 1588.1666
$ws.Range("I132").Value = 1255.5769
$ws.Range("K132").Value = 3766.7307
$ws.Range("M132").Value = -1236.7307

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = $null

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2494.7073
$ws.Range("J132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("N132").Value = $null

$ws.Range("H136").Value = 1714.381
$ws.Range("I136").Value = 1655.6666
$ws.Range("J136").Value = 2066.6667
$ws.Range("K136").Value = 4966.9998
$ws.Range("L136").Value = 6200.000100000001
$ws.Range("M136").Value = -2416.9998
$ws.Range("N136").Value = -11300.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 334796.5
$ws.Range("I81").Value = 334633.34
$ws.Range("K81").Value = 669266.6800000001
$ws.Range("M81").Value = -668205.6800000001

$ws.Range("H84").Value = 334796.5
$ws.Range("I84").Value = 334633.34
$ws.Range("K84").Value = 3346333.4
$ws.Range("M84").Value = -3341029.4

$ws.Range("H132").Value = 4643.6
$ws.Range("I132").Value = 4589.5
$ws.Range("J132").Value = 4860
$ws.Range("K132").Value = 13768.5
$ws.Range("L132").Value = 14580
$ws.Range("M132").Value = -11238.5
$ws.Range("N132").Value = -19640
